$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.202.39'
$ws.Range("E2").Value = '  -0.30%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.681.93'
$ws.Range("E3").Value = '  +0.29%  '

# Row 4
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.95'
$ws.Range("E5").Value = '  -0.64%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5262'
$ws.Range("E6").Value = '  -0.45%  '

# Row 7
$ws.Range("E7").Value = '  -0.03%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2695'
$ws.Range("E8").Value = '  +0.31%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06372'
$ws.Range("E9").Value = '  -1.48%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.52'
$ws.Range("E10").Value = '  -1.66%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07646'
$ws.Range("E11").Value = '  +1.80%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.680.41'
$ws.Range("E12").Value = '  -0.18%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.511'
$ws.Range("E13").Value = '  -0.05%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5754'

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000008331'
$ws.Range("E15").Value = '  -1.97%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.05'
$ws.Range("E16").Value = '  +2.21%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.255.27'
$ws.Range("E17").Value = '  -0.20%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.860'
$ws.Range("E19").Value = '  -1.03%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.81'
$ws.Range("E20").Value = '  -0.34%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '189.09'
$ws.Range("E21").Value = '  -0.40%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.227'
$ws.Range("E22").Value = '  +0.72%  '

# Row 23
$ws.Range("E23").Value = '  -0.01%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '148.74'
$ws.Range("E24").Value = '  +2.71%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.786'
$ws.Range("E25").Value = '  -0.19%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1259'
$ws.Range("E26").Value = '  -1.08%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.72'
$ws.Range("E27").Value = '  -0.17%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06281'
$ws.Range("E28").Value = '  -3.34%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.374'
$ws.Range("E29").Value = '  +0.76%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.316'
$ws.Range("E30").Value = '  -0.15%  '

# Row 31
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.560'
$ws.Range("E31").Value = '  -0.72%  '

# Row 32
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.567'
$ws.Range("E32").Value = '  -0.55%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.682'
$ws.Range("E33").Value = '  +1.73%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.022'
$ws.Range("E34").Value = '  -0.71%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6121'
$ws.Range("E35").Value = '  -1.09%  '

# Row 36
$ws.Range("E36").Value = '  +0.64%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.755'
$ws.Range("E37").Value = '  +0.93%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.181'
$ws.Range("E38").Value = '  -1.57%  '

# Row 39
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01619'
$ws.Range("E39").Value = '  -0.15%  '

# Row 40
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8913'
$ws.Range("E40").Value = '  +2.22%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.098.76'

# Row 42
$ws.Range("E42").Value = '  -0.38%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.41'
$ws.Range("E43").Value = '  -0.03%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.832.79'
$ws.Range("E44").Value = '  +0.37%  '

# Row 45
$ws.Range("E45").Value = '  -0.35%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '57.32'
$ws.Range("E46").Value = '  +0.70%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.004'
$ws.Range("E47").Value = '  -0.37%  '

# Row 48
$ws.Range("E48").Value = '  -1.17%  '

# Row 49
$ws.Range("E49").Value = '  +0.26%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4287'
$ws.Range("E50").Value = '  -0.11%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.012'
$ws.Range("E51").Value = '  -0.69%  '
